$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 893
$ws.Range("I6").Value = 892.125
$ws.Range("J6").Value = 900
$ws.Range("K6").Value = 2676.375
$ws.Range("L6").Value = 2700
$ws.Range("M6").Value = -2564.375
$ws.Range("N6").Value = -2924
$ws.Range("H8").Value = 19844.941
$ws.Range("I8").Value = 66785.07000000001
$ws.Range("J8").Value = 286.55554
$ws.Range("K8").Value = 200355.21
$ws.Range("L8").Value = 859.66662
$ws.Range("M8").Value = -200216.21
$ws.Range("N8").Value = -1137.66662
$ws.Range("H31").Value = 407.83334
$ws.Range("I31").Value = 348.25
$ws.Range("J31").Value = 527
$ws.Range("K31").Value = 1044.75
$ws.Range("L31").Value = 1581
$ws.Range("M31").Value = -814.75
$ws.Range("N31").Value = -2041
$ws.Range("H32").Value = 8598.4
$ws.Range("J32").Value = 8748.25
$ws.Range("L32").Value = 8748.25
$ws.Range("N32").Value = -9400.25
$ws.Range("H33").Value = 1375.25
$ws.Range("I33").Value = 706.2222
$ws.Range("J33").Value = 3382.3333
$ws.Range("K33").Value = 706.2222
$ws.Range("L33").Value = 3382.3333
$ws.Range("M33").Value = -477.2222
$ws.Range("N33").Value = -3840.3333
$ws.Range("H43").Value = 3816.5833
$ws.Range("I43").Value = 3679.9
$ws.Range("J43").Value = 4500
$ws.Range("K43").Value = 3679.9
$ws.Range("L43").Value = 4500
$ws.Range("M43").Value = -3610.9
$ws.Range("N43").Value = -4638
$ws.Range("H99").Value = 1425.1428
$ws.Range("J99").Value = 2933.3333
$ws.Range("L99").Value = 8799.999899999999
$ws.Range("N99").Value = -11795.9999
$ws.Range("H132").Value = 1450917.4
$ws.Range("I132").Value = 1755439.8
$ws.Range("K132").Value = 5266319.4
$ws.Range("M132").Value = -5263789.4
$ws.Range("H137").Value = 3827.359
$ws.Range("I137").Value = 2028.76
$ws.Range("J137").Value = 7039.143
$ws.Range("K137").Value = 6086.28
$ws.Range("L137").Value = 21117.429
$ws.Range("M137").Value = -3536.28
$ws.Range("N137").Value = -26217.429
$ws.Range("H138").Value = 4043.2278
$ws.Range("J138").Value = 5938.41
$ws.Range("L138").Value = 17815.23
$ws.Range("N138").Value = -28095.23

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 89.125
$ws.Range("I5").Value = 101.57143
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 101.57143
$ws.Range("L5").Value = 2
$ws.Range("M5").Value = 10.42856999999999
$ws.Range("N5").Value = -226
$ws.Range("H16").Value = 517
$ws.Range("I16").Value = 338
$ws.Range("J16").Value = 696
$ws.Range("K16").Value = 338
$ws.Range("L16").Value = 696
$ws.Range("M16").Value = -51
$ws.Range("N16").Value = -1270
$ws.Range("H30").Value = 1336.6666
$ws.Range("J30").Value = 2010
$ws.Range("L30").Value = 2010
$ws.Range("N30").Value = -2310
$ws.Range("H45").Value = 3763.923
$ws.Range("I45").Value = 1881.6428
$ws.Range("K45").Value = 1881.6428
$ws.Range("M45").Value = -1504.6428
$ws.Range("H61").Value = 4073.8572
$ws.Range("I61").Value = 3033.875
$ws.Range("J61").Value = 6342.909
$ws.Range("K61").Value = 3033.875
$ws.Range("L61").Value = 6342.909
$ws.Range("M61").Value = -2821.875
$ws.Range("N61").Value = -6766.909
$ws.Range("H97").Value = 1750.3334
$ws.Range("I97").Value = 2000
$ws.Range("J97").Value = 801.6
$ws.Range("K97").Value = 2000
$ws.Range("L97").Value = 801.6
$ws.Range("M97").Value = -1504
$ws.Range("N97").Value = -1793.6
$ws.Range("H132").Value = 4273.2456
$ws.Range("I132").Value = 1759.7693
$ws.Range("K132").Value = 5279.3079
$ws.Range("M132").Value = -2749.3079
$ws.Range("H136").Value = 4073.8572
$ws.Range("I136").Value = 3033.875
$ws.Range("J136").Value = 6342.909
$ws.Range("K136").Value = 9101.625
$ws.Range("L136").Value = 19028.727
$ws.Range("M136").Value = -6551.625
$ws.Range("N136").Value = -24128.727

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 89.125
$ws.Range("I4").Value = 101.57143
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 101.57143
$ws.Range("L4").Value = 2
$ws.Range("M4").Value = 13.42856999999999
$ws.Range("N4").Value = -232
$ws.Range("H82").Value = 40769.25
$ws.Range("H85").Value = 40769.25
$ws.Range("H86").Value = 609276.25
$ws.Range("I86").Value = 1418601.1
$ws.Range("J86").Value = 2282.625
$ws.Range("K86").Value = 1418601.1
$ws.Range("L86").Value = 2282.625
$ws.Range("M86").Value = -1417478.1
$ws.Range("N86").Value = -4528.625
$ws.Range("H89").Value = 609276.25
$ws.Range("I89").Value = 1418601.1
$ws.Range("J89").Value = 2282.625
$ws.Range("K89").Value = 7093005.5
$ws.Range("L89").Value = 11413.125
$ws.Range("M89").Value = -7087389.5
$ws.Range("N89").Value = -22645.125
$ws.Range("H92").Value = 75000
$ws.Range("J92").Value = 75000
$ws.Range("L92").Value = 75000
$ws.Range("N92").Value = -79992
$ws.Range("H94").Value = 6990.1113
$ws.Range("I94").Value = 2150.5
$ws.Range("K94").Value = 2150.5
$ws.Range("M94").Value = -1699.5
$ws.Range("H105").Value = 1583
$ws.Range("I105").Value = 1499.6
$ws.Range("K105").Value = 1499.6
$ws.Range("M105").Value = 247.4000000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 339.6
$ws.Range("J22").Value = 499
$ws.Range("L22").Value = 499
$ws.Range("N22").Value = -1199
$ws.Range("H44").Value = 2000
$ws.Range("I44").Value = 2000
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 2000
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -1558
$ws.Range("N44").ClearContents()
$ws.Range("H94").Value = 1026.3914
$ws.Range("I94").Value = 1268.5555
$ws.Range("J94").Value = 870.7143
$ws.Range("K94").Value = 1268.5555
$ws.Range("L94").Value = 870.7143
$ws.Range("M94").Value = -817.5554999999999
$ws.Range("N94").Value = -1772.7143
$ws.Range("H105").Value = 1947.6666
$ws.Range("I105").Value = 1421.5
$ws.Range("K105").Value = 1421.5
$ws.Range("M105").Value = 325.5
$ws.Range("H139").Value = 95997.25
$ws.Range("I139").Value = 95000
$ws.Range("J139").Value = 96994.5
$ws.Range("K139").Value = 95000
$ws.Range("L139").Value = 96994.5
$ws.Range("M139").Value = -89860
$ws.Range("N139").Value = -107274.5
$ws.Range("H141").Value = 87373
$ws.Range("J141").Value = 96283.42999999999
$ws.Range("L141").Value = 96283.42999999999
$ws.Range("N141").Value = -106643.43

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 60310
$ws.Range("J34").Value = 92692.30499999999
$ws.Range("L34").Value = 278076.915
$ws.Range("N34").Value = -278244.915
$ws.Range("H41").Value = 1167.5
$ws.Range("I41").Value = 335
$ws.Range("J41").Value = 2000
$ws.Range("K41").Value = 1005
$ws.Range("L41").Value = 6000
$ws.Range("M41").Value = -667
$ws.Range("N41").Value = -6676
$ws.Range("H103").Value = 242.7
$ws.Range("I103").Value = 182
$ws.Range("K103").Value = 546
$ws.Range("M103").Value = 333
$ws.Range("H107").Value = 48429.273
$ws.Range("I107").Value = 1452.3334
$ws.Range("J107").Value = 66045.625
$ws.Range("K107").Value = 4357.0002
$ws.Range("L107").Value = 198136.875
$ws.Range("M107").Value = -2437.0002
$ws.Range("N107").Value = -201976.875
$ws.Range("H134").Value = 1196.7693
$ws.Range("I134").Value = 888.1667
$ws.Range("J134").Value = 4900
$ws.Range("K134").Value = 2664.5001
$ws.Range("L134").Value = 14700
$ws.Range("M134").Value = 2405.4999
$ws.Range("N134").Value = -24840
$ws.Range("H136").Value = 8924.5
$ws.Range("I136").Value = 8899.333000000001
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 26697.999
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -21597.999
$ws.Range("N136").Value = -37200

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2186.9375
$ws.Range("I102").Value = 1779.3
$ws.Range("K102").Value = 1779.3
$ws.Range("M102").Value = -157.3
$ws.Range("H122").Value = 4449
$ws.Range("I122").Value = 3869.9285
$ws.Range("J122").Value = 5607.143
$ws.Range("K122").Value = 11609.7855
$ws.Range("L122").Value = 16821.429
$ws.Range("M122").Value = -9159.7855
$ws.Range("N122").Value = -21721.429
$ws.Range("H132").Value = 442415.44
$ws.Range("I132").Value = 478646.28
$ws.Range("J132").Value = 252203.5
$ws.Range("K132").Value = 1435938.84
$ws.Range("L132").Value = 756610.5
$ws.Range("M132").Value = -1433408.84
$ws.Range("N132").Value = -761670.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 81492.53999999999
$ws.Range("I40").Value = 170336
$ws.Range("J40").Value = 5341
$ws.Range("K40").Value = 170336
$ws.Range("L40").Value = 5341
$ws.Range("M40").Value = -170200
$ws.Range("N40").Value = -5613
$ws.Range("H100").Value = 1000
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 1000000000
$ws.Range("J29").Value = 1000000000
$ws.Range("L29").Value = 1000000000
$ws.Range("N29").Value = -1000000580
$ws.Range("H100").Value = 462.52942
$ws.Range("I100").Value = 418.72726
$ws.Range("K100").Value = 837.45452
$ws.Range("M100").Value = -296.45452
$ws.Range("H133").Value = 62528.57
$ws.Range("J133").Value = 62528.57
$ws.Range("L133").Value = 62528.57
$ws.Range("N133").Value = -72648.57000000001
